$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Force text type (matches original inlineStr cells) even for values that
# look numeric, then reset the style so no new "s" attribute is introduced.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.101.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.891.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9989'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5183'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3729'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07204'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9060'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.09'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07657'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.892.02'
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.277'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9991'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008505'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.149.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.059'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.106.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.470'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.791'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.152'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.938'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.809'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09210'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05051'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7628'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.196'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.023'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.286'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.560'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5625'
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.599'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '118.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.884'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1507'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4804'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9986'
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.72'
$ws.Range("D51").Style = "Normal"

# --- Volume / 1h change (column E) updates ---
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("E7").Value = '  +2.10%  '
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("E10").Value = '  +1.87%  '
$ws.Range("E11").Value = '  +1.98%  '
$ws.Range("E12").Value = '  +1.66%  '
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("E14").Value = '  +3.61%  '
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("E18").Value = '  +1.96%  '
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("E21").Value = '  +0.92%  '
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("E23").Value = '  +2.47%  '
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  -0.79%  '
$ws.Range("E28").Value = '  +5.10%  '
$ws.Range("E29").Value = '  +1.45%  '
$ws.Range("E30").Value = '  +5.76%  '
$ws.Range("E31").Value = '  +3.64%  '
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("E33").Value = '  -0.86%  '
$ws.Range("E34").Value = '  +3.97%  '
$ws.Range("E35").Value = '  +4.07%  '
$ws.Range("E36").Value = '  -2.30%  '
$ws.Range("E37").Value = '  +2.29%  '
$ws.Range("E38").Value = '  +4.06%  '
$ws.Range("E39").Value = '  +5.85%  '
$ws.Range("E40").Value = '  -0.62%  '
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("E44").Value = '  +5.92%  '
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("E46").Value = '  +3.67%  '
$ws.Range("E47").Value = '  +2.26%  '
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("E49").Value = '  +1.22%  '
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("E51").Value = '  +1.31%  '
